$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ขนส่ง" (AA) and "Tracking" (AB) columns from the export template.
# Deleting the entire columns shifts the remaining headers (old AC/AD, i.e.
# "สถานะการชำระเงิน" / "วันที่รับชำระเงิน") left into AA/AB and drops the
# now-unused shared strings automatically.
$ws.Range("AA1:AB1").EntireColumn.Delete()

# Leave the sheet's recorded selection on A2, matching the saved workbook state.
$ws.Range("A2").Select()
